$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "28.260.15"
$ws.Range("E2").Value = "  -0.52%  "
# Row 3
$ws.Range("D3").Value = "1.805.79"
$ws.Range("E3").Value = "  -0.82%  "
# Row 4
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "
# Row 5
$ws.Range("D5").Value = "'314.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.28%  "
# Row 6
$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.01%  "
# Row 7
$ws.Range("D7").Value = "'0.5274"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.05%  "
# Row 8
$ws.Range("D8").Value = "'0.3843"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.03%  "
# Row 9
$ws.Range("D9").Value = "'0.08038"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.28%  "
# Row 10
$ws.Range("D10").Value = "'41.53"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.57%  "
# Row 11
$ws.Range("D11").Value = "'1.105"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.36%  "
# Row 12
$ws.Range("D12").Value = "'6.355"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.75%  "
# Row 13
$ws.Range("D13").Value = "'1.001"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.03%  "
# Row 14
$ws.Range("D14").Value = "'20.68"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.21%  "
# Row 15
$ws.Range("D15").Value = "'7.354"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.68%  "
# Row 16
$ws.Range("D16").Value = "1.803.33"
$ws.Range("E16").Value = "  -1.10%  "
# Row 17
$ws.Range("D17").Value = "'92.89"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.40%  "
# Row 18
$ws.Range("D18").Value = "'0.00001098"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.30%  "
# Row 19
$ws.Range("D19").Value = "'0.06609"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.17%  "
# Row 20
$ws.Range("D20").Value = "'1.001"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.05%  "
# Row 21
$ws.Range("D21").Value = "'17.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.42%  "
# Row 22
$ws.Range("D22").Value = "'5.991"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.48%  "
# Row 23
$ws.Range("D23").Value = "28.293.39"
$ws.Range("E23").Value = "  -0.46%  "
# Row 24
$ws.Range("D24").Value = "'11.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.19%  "
# Row 25
$ws.Range("D25").Value = "'2.234"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.10%  "
# Row 26
$ws.Range("D26").Value = "'160.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.54%  "
# Row 27
$ws.Range("D27").Value = "'20.58"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.54%  "
# Row 28
$ws.Range("D28").Value = "2.008.15"
$ws.Range("E28").Value = "  -0.80%  "
# Row 29
$ws.Range("D29").Value = "'2.396"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.06%  "
# Row 30
$ws.Range("D30").Value = "'123.36"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.59%  "
# Row 31
$ws.Range("D31").Value = "'0.1088"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.13%  "
# Row 32
$ws.Range("D32").Value = "'1.063"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.46%  "
# Row 33
$ws.Range("D33").Value = "'3.668"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.52%  "
# Row 34
$ws.Range("D34").Value = "'5.585"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.50%  "
# Row 35
$ws.Range("D35").Value = "'0.07284"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.18%  "
# Row 36
$ws.Range("D36").Value = "'12.33"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.80%  "
# Row 37
$ws.Range("B37").Value = "Algorand"
$ws.Range("C37").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D37").Value = "'0.2173"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.85%  "
# Row 38
$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").Value = "'8.876"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.27%  "
# Row 39
$ws.Range("D39").Value = "'0.02320"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.06%  "
# Row 40
$ws.Range("D40").Value = "'5.117"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.25%  "
# Row 41
$ws.Range("D41").Value = "'0.6235"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.32%  "
# Row 42
$ws.Range("D42").Value = "'1.168"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.51%  "
# Row 43
$ws.Range("D43").Value = "'1.370"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.54%  "
# Row 44
$ws.Range("D44").Value = "'13.28"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.40%  "
# Row 45
$ws.Range("D45").Value = "'0.6028"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.51%  "
# Row 46
$ws.Range("D46").Value = "'3.768"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.95%  "
# Row 47
$ws.Range("D47").Value = "'127.07"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.29%  "
# Row 48
$ws.Range("D48").Value = "'1.215"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.76%  "
# Row 49
$ws.Range("D49").Value = "'1.936"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.93%  "
# Row 50
$ws.Range("D50").Value = "'0.06840"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.77%  "
# Row 51
$ws.Range("D51").Value = "'73.18"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.14%  "
